$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8464730381965637
$ws.Range("B1").Value = 1.695471882820129
$ws.Range("C1").Value = 6.239840507507324
$ws.Range("D1").Value = 1.901181221008301
$ws.Range("E1").Value = 1.152170896530151
